$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(15).Delete()
$ws.Range("E15").Select()
